# Update for release to deploy 0.1.1
$wb = $excel.ActiveWorkbook

# --- Rename the "Include ValueSets" sheets ---
$wb.Worksheets.Item("Include ValueSets").Name = "Include ValueSet #0"
$wb.Worksheets.Item("Include ValueSets 2").Name = "Include ValueSet #1"
$wb.Worksheets.Item("Include ValueSets 3").Name = "Include ValueSet #2"

# --- Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Bump version (row 3: Version / 0.1.0 -> 0.1.1)
$ws.Range("B3").Value = "0.1.1"

# Bump date (row 8: Date / ... -> 2024-11-11T17:53:38-06:00)
$ws.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new "Jurisdiction" row (with blank value) right after the "Contact"
# row (row 10), pushing Description/Purpose/Copyright/Immutable down by one.
$ws.Rows.Item(11).Insert()

# Match the look of the surrounding property rows as closely as possible.
$newRow = $ws.Range("A11:B11")
$newRow.VerticalAlignment = -4160
$newRow.WrapText = $true
$newRow.Borders.LineStyle = 1
$newRow.Borders.ColorIndex = 23

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# The rest of the rows below keep their original content - only their row
# numbers shifted down by the insert above.
